$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting
# (avoids Excel auto-converting numeric-looking strings to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "58.523.47"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "2.611.71"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "534.09"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "142.54"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").Value = "6.85"
$ws.Range("E9").Value = "  +4.85%  "
$ws.Range("D10").Value = "0.0997"
$ws.Range("E10").Value = "  -2.67%  "
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").Value = "3.080.79"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "58.461.35"
$ws.Range("D15").Value = "20.69"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "2.622.98"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").Value = "333.71"
$ws.Range("E19").Value = "  -2.22%  "
$ws.Range("D20").Value = "10.10"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "66.42"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("E26").Value = "  -2.28%  "
$ws.Range("D27").Value = "7.07"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D29").Value = "0.0₃0731"
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").Value = "5.85"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "152.87"
$ws.Range("E32").Value = "  +2.08%  "
$ws.Range("D33").Value = "18.84"
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("D34").Value = "3.88"
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("D36").Value = "0.839"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "1.41"
$ws.Range("E37").Value = "  -2.99%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "0.812"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").Value = "3.56"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").Value = "281.53"
$ws.Range("E40").Value = "  +3.11%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "0.593"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").Value = "0.0940"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D45").Value = "18.92"
$ws.Range("E45").Value = "  +1.87%  "
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("D48").Value = "1.939.90"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "4.43"
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("D50").Value = "17.78"
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("D51").Value = "113.72"
$ws.Range("E51").Value = "  +0.57%  "
